$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Apply formatting (copy existing style definitions onto target ranges) ---
# Source cells are stable originals that keep their style: D1 (style 4), J1 (style 1),
# D4 (style 2), D7 (style 3). Copying format only (no values) reuses existing style/font
# table entries instead of creating new ones.

$ws.Range("D1").Copy() | Out-Null
$ws.Range("D1:I1").PasteSpecial(-4122) | Out-Null

$ws.Range("J1").Copy() | Out-Null
$ws.Range("J1:N1").PasteSpecial(-4122) | Out-Null

$ws.Range("D4").Copy() | Out-Null
$ws.Range("B3:J3").PasteSpecial(-4122) | Out-Null
$ws.Range("D4:J4").PasteSpecial(-4122) | Out-Null
$ws.Range("D5:E5").PasteSpecial(-4122) | Out-Null
$ws.Range("G5:J5").PasteSpecial(-4122) | Out-Null
$ws.Range("E6:J6").PasteSpecial(-4122) | Out-Null
$ws.Range("F7:J7").PasteSpecial(-4122) | Out-Null
$ws.Range("E8:J8").PasteSpecial(-4122) | Out-Null
$ws.Range("E9:J9").PasteSpecial(-4122) | Out-Null
$ws.Range("E10:J10").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("F11").PasteSpecial(-4122) | Out-Null
$ws.Range("F12").PasteSpecial(-4122) | Out-Null
$ws.Range("F13").PasteSpecial(-4122) | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null
$ws.Range("B43").PasteSpecial(-4122) | Out-Null

$ws.Range("D7").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("D7:E7").PasteSpecial(-4122) | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Step 2: Set cell text values ---
$ws.Range("D1").Value = "CATECHISM OF THE CATHOLIC CHURCH"
$ws.Range("B3").Value = "ADMIN"
$ws.Range("D4").Value = "Login"
$ws.Range("D5").Value = "Granting the request for teachers registration"
$ws.Range("D6").Value = "Manage Examination Details"
$ws.Range("F7").Value = "Send Examination Notifications to Teachers"
$ws.Range("F8").Value = "Send Examination Notifications to Parents"
$ws.Range("D9").Value = "Manage Extra Carricular Activities"
$ws.Range("F10").Value = "Publishing Points of Activities Based on Groups"
$ws.Range("D11").Value = "Notifying Details of Each Programs"
$ws.Range("F12").Value = "send Notifications to Teachers"
$ws.Range("D13").Value = "Rank List Publication"
$ws.Range("F14").Value = "Result for All Classes"
$ws.Range("D15").Value = "Leave Confirmation"
$ws.Range("F16").Value = "Approve Leave Rquest from Teachers"
$ws.Range("D17").Value = "Manage Details of Group Division"
$ws.Range("D18").Value = "Manage Details of Mission League"
$ws.Range("D19").Value = "Manage Details of Thirubalasakhyam"
$ws.Range("B21").Value = "STAFF"
$ws.Range("D22").Value = "Login"
$ws.Range("D23").Value = "Add Students"
$ws.Range("D24").Value = "Students Information"
$ws.Range("D25").Value = "Report Generation"
$ws.Range("F26").Value = "Generate Reports"
$ws.Range("F27").Value = "Send Reports to Admin"
$ws.Range("D28").Value = "Attendance Monitoring"
$ws.Range("F29").Value = "Students Attendance"
$ws.Range("D30").Value = "Leave Requisition"
$ws.Range("F31").Value = "To Admin"
$ws.Range("D32").Value = "Leave Confirmation"
$ws.Range("F33").Value = "From Parents"
$ws.Range("D34").Value = "View Rank List"
$ws.Range("F35").Value = "Result for All Classes"
$ws.Range("D36").Value = "Notifications"
$ws.Range("F37").Value = "View Notifications From Admin"
$ws.Range("F38").Value = "Send Notifications to Parents"
$ws.Range("D39").Value = "View Examination Details"
$ws.Range("D40").Value = "Group Division"
$ws.Range("F41").Value = "View Group Division Lists"
$ws.Range("B43").Value = "PARENT"
$ws.Range("D44").Value = "Login"
$ws.Range("D45").Value = "View Rank List"
$ws.Range("D46").Value = "View Reports"
$ws.Range("D47").Value = "Leave Requisition"
$ws.Range("F48").Value = "To Teachers"
$ws.Range("D49").Value = "View Examination Details"
$ws.Range("F50").Value = "View Dates of Exams"
$ws.Range("D51").Value = "Parent Information"

# --- Step 3: Cells that must be blank but keep their (already-correct) style ---
$ws.Range("E1").ClearContents()
$ws.Range("F1").ClearContents()
$ws.Range("G1").ClearContents()
$ws.Range("H1").ClearContents()
$ws.Range("I1").ClearContents()
$ws.Range("J1").ClearContents()
$ws.Range("K1").ClearContents()
$ws.Range("L1").ClearContents()
$ws.Range("M1").ClearContents()
$ws.Range("N1").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("G10").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("I10").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("B11").ClearContents()
$ws.Range("F11").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("D16").ClearContents()

# --- Step 4: Fully clear cells orphaned by the restructure (no longer used at all) ---
$ws.Range("F5").Clear()
$ws.Range("F17").Clear()
$ws.Range("F18").Clear()
$ws.Range("F20").Clear()
$ws.Range("F21").Clear()
$ws.Range("F23").Clear()
$ws.Range("F24").Clear()
$ws.Range("D26").Clear()
$ws.Range("D27").Clear()
$ws.Range("B28").Clear()
$ws.Range("D29").Clear()
$ws.Range("F30").Clear()
$ws.Range("D31").Clear()

# --- Step 5: Row heights for newly (re)created rows that need the 15.75 pt height ---
$ws.Rows(18).RowHeight = 15.75
$ws.Rows(21).RowHeight = 15.75
$ws.Rows(22).RowHeight = 15.75
$ws.Rows(43).RowHeight = 15.75

# --- Step 6: Update selection to match the author's final cursor position ---
$ws.Range("B43").Select() | Out-Null
